$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

Set-TextValue $ws.Range('D2') '51.367.65'
$ws.Range('E2').Value = '  -0.66%  '
Set-TextValue $ws.Range('D3') '3.092.20'
$ws.Range('E3').Value = '  +2.05%  '
$ws.Range('E4').Value = '  +0.06%  '
Set-TextValue $ws.Range('D5') '386.55'
$ws.Range('E5').Value = '  +1.62%  '
Set-TextValue $ws.Range('D6') '103.45'
$ws.Range('E6').Value = '  +0.55%  '
Set-TextValue $ws.Range('D7') '0.537'
$ws.Range('E7').Value = '  -1.60%  '
$ws.Range('E8').Value = '  +0.04%  '
Set-TextValue $ws.Range('D9') '0.585'
$ws.Range('E9').Value = '  -1.51%  '
Set-TextValue $ws.Range('D10') '36.89'
$ws.Range('E10').Value = '  +0.15%  '
$ws.Range('E11').Value = '  +0.04%  '
Set-TextValue $ws.Range('D12') '0.0855'
$ws.Range('E12').Value = '  -0.54%  '
Set-TextValue $ws.Range('D13') '3.586.14'
$ws.Range('E13').Value = '  +2.12%  '
Set-TextValue $ws.Range('D14') '18.47'
$ws.Range('E14').Value = '  -0.48%  '
Set-TextValue $ws.Range('D15') '7.76'
$ws.Range('E15').Value = '  +0.26%  '
Set-TextValue $ws.Range('D16') '3.099.35'
$ws.Range('E16').Value = '  +2.26%  '
Set-TextValue $ws.Range('D17') '0.991'
$ws.Range('E17').Value = '  +1.40%  '
Set-TextValue $ws.Range('D18') '10.63'
$ws.Range('E18').Value = '  +0.65%  '
Set-TextValue $ws.Range('D19') '51.467.58'
$ws.Range('E19').Value = '  -0.46%  '
Set-TextValue $ws.Range('D20') '3.24'
$ws.Range('E20').Value = '  +5.54%  '
Set-TextValue $ws.Range('D21') '12.45'
$ws.Range('E21').Value = '  -0.47%  '
Set-TextValue $ws.Range('D22') '0.0₃0963'
$ws.Range('E22').Value = '  +0.12%  '
Set-TextValue $ws.Range('D23') '70.04'
$ws.Range('E23').Value = '  -0.10%  '
Set-TextValue $ws.Range('D24') '265.66'
$ws.Range('E24').Value = '  -1.05%  '
Set-TextValue $ws.Range('D25') '3.15'
$ws.Range('E25').Value = '  -0.18%  '
Set-TextValue $ws.Range('D26') '7.97'
$ws.Range('E26').Value = '  -3.52%  '
Set-TextValue $ws.Range('D27') '27.26'
$ws.Range('E27').Value = '  +3.82%  '
Set-TextValue $ws.Range('D28') '0.999'
$ws.Range('E28').Value = '  -0.03%  '
Set-TextValue $ws.Range('D29') '7.15'
$ws.Range('E29').Value = '  -6.50%  '
$ws.Range('E30').Value = '  -4.93%  '
Set-TextValue $ws.Range('D31') '0.106'
$ws.Range('E31').Value = '  -2.19%  '
Set-TextValue $ws.Range('D32') '10.38'
$ws.Range('E32').Value = '  +0.81%  '
Set-TextValue $ws.Range('D33') '35.53'
$ws.Range('E33').Value = '  +4.08%  '
Set-TextValue $ws.Range('D34') '0.0471'
$ws.Range('E34').Value = '  +5.11%  '
$ws.Range('E35').Value = '  +1.87%  '
Set-TextValue $ws.Range('D36') '49.97'
$ws.Range('E36').Value = '  -1.17%  '
$ws.Range('E37').Value = '  -0.03%  '
Set-TextValue $ws.Range('D38') '3.36'
$ws.Range('E38').Value = '  +1.40%  '
Set-TextValue $ws.Range('D39') '0.289'
$ws.Range('E39').Value = '  -1.51%  '
Set-TextValue $ws.Range('D40') '129.40'
$ws.Range('E40').Value = '  +4.43%  '
Set-TextValue $ws.Range('D41') '1.85'
$ws.Range('E41').Value = '  -0.65%  '
Set-TextValue $ws.Range('D42') '0.115'
$ws.Range('E42').Value = '  -0.39%  '
Set-TextValue $ws.Range('D43') '16.50'
$ws.Range('E43').Value = '  -3.25%  '
Set-TextValue $ws.Range('D44') '3.80'
$ws.Range('E44').Value = '  +0.93%  '
Set-TextValue $ws.Range('D45') '2.49'
$ws.Range('E45').Value = '  -2.99%  '
Set-TextValue $ws.Range('D46') '22.01'
$ws.Range('E46').Value = '  +0.84%  '
$ws.Range('E47').Value = '  +3.70%  '
$ws.Range('E48').Value = '  -0.36%  '
Set-TextValue $ws.Range('D49') '2.072.49'
$ws.Range('E49').Value = '  +1.84%  '
$ws.Range('E50').Value = '  +3.95%  '
$ws.Range('E51').Value = '  +18.12%  '
